# Applies the "added personal information" edit to the Personal Information document.
# Samuel Everson's bio paragraph is re-split into multiple runs (matching Word's
# proofing-mark behaviour) with a "self taught" -> "self-taught" correction, and a
# new blank paragraph plus Stanton Wightwick's name/student-number line and full
# bio paragraph are appended (carrying the "_GoBack" bookmark to its new position).

$d = $word.ActiveDocument

$p5 = $d.Paragraphs(5)
$r = $p5.Range

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">I’m 27, live on the northern end of the Gold Coast. I have 2 daughters and a spouse. I like tinkering with things like motors, computers and various other gadgets (though I’m still learning a lot and am not very </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>good</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>!). I enjoy the idea of software engineering though have very little experience in it. In the same boat I also like database design and implementation. Unfortunately other than a Diploma in IT General I have no formal</w:t></w:r><w:r><w:t xml:space="preserve"> experience in IT but have self-</w:t></w:r><w:r><w:t>taught a range of IT skills from building desktop PC’s and servers, VBA, HTML/CSS and database integration using MySQL.  I’m very happy to be a part of &lt;insert team name&gt; and am looking forward to working on this project with the team!</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Stanton </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Wightwick</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – S3819611</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I’m 26, hail from Melbourne CBD. I am a huge basketball and combat sports fan, having played the former most my life and also competed in the latter a couple of times. </w:t></w:r><w:r><w:t xml:space="preserve">I also have an unhealthy obsession with watches (save your wallet, do not get into this hobby!!) </w:t></w:r><w:r><w:t>My interest in I.T has definitely been amplified through my full time job as a banker, and seeing the industry’s focus shift to improvement in technology,</w:t></w:r><w:r><w:t xml:space="preserve"> systems etc.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Currently, (</w:t></w:r><w:r><w:t>which i</w:t></w:r><w:r><w:t xml:space="preserve">ntro to RMIT has helped me realise) I enjoy more so the front end development/software dev side. </w:t></w:r><w:r><w:t>I am drawn to</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> ‘ease of life’ type applications and nicely designed, interactive interfaces. I have no formal education, my minimal I.T skills are mostly self-taught (Excel, VBA, HTML/CSS, Python) either job necessity or through websites such as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Udemy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, EDX, YouTube etc. I’m excited to work with &lt;insert team name&gt;, and can’t wait to see what we come up with!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xml)
